$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the B-column weight labels so they follow the new assignment pattern
$ws.Range("B2").Value = "w_1"
$ws.Range("B3").Value = "w_1"
$ws.Range("B4").Value = "w_1"
$ws.Range("B5").Value = "w_1"
$ws.Range("B9").Value = "w_2"
$ws.Range("B14").Value = "w_3"
$ws.Range("B19").Value = "w_3"
$ws.Range("B27").Value = "w_4"

# Remove the left indent on the weight column so any value greater than 1 fits
$ws.Range("B2:B31").IndentLevel = 0

# Reflect where the user last clicked
$ws.Range("B27").Select()
